# Applies the "Added SysID transformer models and transient operating mode" commit.
$wb = $excel.ActiveWorkbook

$nl = [char]10

# ---------------------------------------------------------------------------
# Exp sheet
# ---------------------------------------------------------------------------
$wsExp = $wb.Worksheets.Item("Exp")

# Waveform Smoothing now defaults to fundamental-cycle averaging.
$wsExp.Range("E4").Value = "fel"

# Simulation Frequency upper bound raised 10M -> 1G.
$wsExp.Range("E5").Value = 1000000000
$wsExp.Range("E5").Validation.Modify(1, 1, 1, 10000, 100000000000)

# New row: "Plot Downsampling" toggle, formatted like the row above it.
$wsExp.Range("A8:F8").Copy()
$wsExp.Range("A9:F9").PasteSpecial(-4122)
$wsExp.Range("A9").Value = "Plot Downsampling"
$wsExp.Range("C9").Value = "0: no downsampling, 1: downsampling"
$wsExp.Range("D9").Value = "bDS"
$wsExp.Range("E9").Value = 0
$wsExp.Range("F9").Value = "-"
$wsExp.Range("E9").Validation.Add(1, 1, 1, 0, 1)

# ---------------------------------------------------------------------------
# Dat sheet
# ---------------------------------------------------------------------------
$wsDat = $wb.Worksheets.Item("Dat")
$wsDat.Range("E6").Value = 2000
$wsDat.Range("E9").Value = 150
$wsDat.Range("E10").Value = 0.5
$wsDat.Range("E11").Value = 425
$wsDat.Range("E13").Value = 0.005
$wsDat.Activate()
$wsDat.Range("C38").Select()

# ---------------------------------------------------------------------------
# Top sheet
# ---------------------------------------------------------------------------
$wsTop = $wb.Worksheets.Item("Top")
$wsTop.Range("E11").Value = 0.05
$wsTop.Range("E12").Value = 0.0000005
$wsTop.Activate()
$wsTop.Range("E2").Select()

# ---------------------------------------------------------------------------
# Par sheet
# ---------------------------------------------------------------------------
$wsPar = $wb.Worksheets.Item("Par")
$wsPar.Range("E12").Value = 100000
$wsPar.Activate()
$wsPar.Range("G2").Select()

# ---------------------------------------------------------------------------
# Mag sheet - new "SS" (state space) transformer model option.
# ---------------------------------------------------------------------------
$wsMag = $wb.Worksheets.Item("Mag")
$wsMag.Range("C2").Value = "(NT): no transformer, " + $nl + "(OC): open circuit at secondary, " + $nl + "(SC): short circuit at secondary, " + $nl + "(RL): RL load at secondary, " + $nl + "(SS): defined via state space model in para-Excel"
$wsMag.Range("E2").Value = "SS"
$wsMag.Range("E2").Validation.Modify(3, 1, 1, '"NT, OC, SC, RL, SS"')

# Mag becomes the active/selected sheet (was Top before).
$wsMag.Activate()
$wsMag.Range("E2").Select()
